$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 321, pushing existing rows 321-328 down to 323-330
$ws.Rows.Item(321).EntireRow.Insert()
$ws.Rows.Item(321).EntireRow.Insert()

# New row 321: Betarraga, "Primera" quality, week of 2022-07-05 (serial 44747)
$ws.Cells.Item(321,1).Value = 1
$ws.Cells.Item(321,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(321,3).Value = "Arica y Parinacota"
$ws.Cells.Item(321,4).Value = 44747
$ws.Cells.Item(321,5).Value = 15
$ws.Cells.Item(321,6).Value = 100114014
$ws.Cells.Item(321,7).Value = "Betarraga"
$ws.Cells.Item(321,8).Value = "Sin especificar"
$ws.Cells.Item(321,9).Value = "Primera"
$ws.Cells.Item(321,10).Value = 1000
$ws.Cells.Item(321,11).Value = 450
$ws.Cells.Item(321,12).Value = 500
$ws.Cells.Item(321,13).Value = 475
$ws.Cells.Item(321,14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(321,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(321,16).Value = 119
$ws.Cells.Item(321,17).Value = 4
$ws.Cells.Item(321,18).Value = "Hortaliza"

# New row 322: Betarraga, "Segunda" quality, same week (serial 44747)
$ws.Cells.Item(322,1).Value = 1
$ws.Cells.Item(322,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(322,3).Value = "Arica y Parinacota"
$ws.Cells.Item(322,4).Value = 44747
$ws.Cells.Item(322,5).Value = 15
$ws.Cells.Item(322,6).Value = 100114014
$ws.Cells.Item(322,7).Value = "Betarraga"
$ws.Cells.Item(322,8).Value = "Sin especificar"
$ws.Cells.Item(322,9).Value = "Segunda"
$ws.Cells.Item(322,10).Value = 900
$ws.Cells.Item(322,11).Value = 450
$ws.Cells.Item(322,12).Value = 500
$ws.Cells.Item(322,13).Value = 475
$ws.Cells.Item(322,14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(322,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(322,16).Value = 95
$ws.Cells.Item(322,17).Value = 5
$ws.Cells.Item(322,18).Value = "Hortaliza"
